$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "67.662.43"
$ws.Range("E2").Value = "  +1.16%  "
$ws.Range("D3").Value = "3.877.10"
$ws.Range("E3").Value = "  +0.82%  "
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "1.00"
$ws.Range("D4").Style = "Normal"
$ws.Range("E4").Value = "  -0.01%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "464.73"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  +9.73%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "148.04"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  +14.10%  "
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.633"
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = "  +4.22%  "
$ws.Range("E8").Value = "  +0.00%  "
$ws.Range("E9").Value = "  +4.43%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "0.157"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "  +0.10%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.0000317"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "  -5.59%  "
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "44.11"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = "  +8.48%  "
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "10.47"
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = "  +2.13%  "
$ws.Range("D14").Value = "4.499.98"
$ws.Range("E14").Value = "  +1.01%  "
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "14.78"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "  -6.30%  "
$ws.Range("D16").Value = "3.883.10"
$ws.Range("E16").Value = "  +0.85%  "
$ws.Range("E17").Value = "  +0.02%  "
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "20.14"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = "  +1.70%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "1.17"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "  +7.85%  "
$ws.Range("D20").Value = "67.780.67"
$ws.Range("E20").Value = "  +0.86%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "432.97"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "  +4.27%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "14.92"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "  -0.08%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "3.29"
$ws.Range("D23").Style = "Normal"
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "88.40"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "  +5.15%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "3.57"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "  +10.48%  "
$ws.Range("B26").Value = "EthereumClassic"
$ws.Range("C26").Value = "https://coinranking.com/coin/hnfQfsYfeIGUQ+ethereumclassic-etc"
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "37.80"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "  +0.98%  "
$ws.Range("B27").Value = "RenderToken"
$ws.Range("C27").Value = "https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr"
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "10.34"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "  +11.90%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "10.27"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = "  +4.37%  "
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "5.52"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = "  +3.77%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "747.39"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "  +1.73%  "
$ws.Range("B31").Value = "Cosmos"
$ws.Range("C31").Value = "https://coinranking.com/coin/Knsels4_Ol-Ny+cosmos-atom"
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "13.88"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "  +6.15%  "
$ws.Range("B32").Value = "Hedera"
$ws.Range("C32").Value = "https://coinranking.com/coin/jad286TjB+hedera-hbar"
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "0.135"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "  +10.76%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "2.76"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "  -0.44%  "
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "43.41"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "  +12.53%  "
$ws.Range("E35").Value = "  +7.62%  "
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "57.49"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = "  +3.67%  "
$ws.Range("B37").Value = "Dai"
$ws.Range("C37").Value = "https://coinranking.com/coin/MoTuySvg7+dai-dai"
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "1.00"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = "  +0.19%  "
$ws.Range("B38").Value = "NEARProtocol"
$ws.Range("C38").Value = "https://coinranking.com/coin/DCrsaMv68+nearprotocol-near"
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "5.62"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "  +2.21%  "
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "0.0481"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "  +4.57%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "0.351"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "  +11.44%  "
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "2.94"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "  +2.22%  "
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "2.61"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "  +14.00%  "
$ws.Range("B43").Value = "Stellar"
$ws.Range("C43").Value = "https://coinranking.com/coin/f3iaFeCKEmkaZ+stellar-xlm"
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "0.142"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "  +5.99%  "
$ws.Range("B44").Value = "PEPE"
$ws.Range("C44").Value = "https://coinranking.com/coin/03WI8NQPF+pepe-pepe"
$ws.Range("D44").Value = "0.0₃0682"
$ws.Range("E44").Value = "  -6.26%  "
$ws.Range("E45").Value = "  +0.06%  "
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "3.30"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "  +6.21%  "
$ws.Range("E47").Value = "  +2.80%  "
$ws.Range("E48").Value = "  +7.76%  "
$ws.Range("E49").Value = "  +4.15%  "
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "145.32"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "  +3.27%  "
$ws.Range("E51").Value = "  +3.86%  "
